# day 1 is complete
# Applies the NACL / Security Groups slide updates (slide 2) and the
# "What are we using today?" SG-target list update (slide 6).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2 - "So, what are they?"
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)

# --- Shape: NACL content placeholder -----------------------------------
$naclShape = $slide2.Shapes.Item(2)
$naclShape.TextFrame.AutoSize = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/>

$naclTr = $naclShape.TextFrame.TextRange
$naclLast = $naclTr.Paragraphs($naclTr.Paragraphs().Count, 1)
$naclLast.InsertAfter("`rGood for controlling OUTBOUND traffic from a locked down subnet i.e. controlling what a server can access") | Out-Null

$naclTr = $naclShape.TextFrame.TextRange
$naclLast = $naclTr.Paragraphs($naclTr.Paragraphs().Count, 1)
$naclLast.InsertAfter("`rDefaults to all ALLOW") | Out-Null

# --- Shape: Security Groups content placeholder -------------------------
$sgShape = $slide2.Shapes.Item(3)
$sgShape.TextFrame.AutoSize = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/>

$sgTr = $sgShape.TextFrame.TextRange
$sgFull = $sgTr.Text
$sgTarget = " you can associate an SG for an EC2"
$sgStart = $sgFull.IndexOf($sgTarget) + 1
$sgSub = $sgTr.Characters($sgStart, $sgTarget.Length)
$sgSub.Text = " you can associate an SG for an EC2 or a Load Balancer"

# ---------------------------------------------------------------------
# Slide 6 - "What are we using today?"
# ---------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$listShape = $slide6.Shapes.Item(4)
$listTr = $listShape.TextFrame.TextRange

# Remove the "SSM Agent IAM Role" bullet entirely.
$listTr.Paragraphs(3, 1).Delete()

# After "Webservers" add "Port 80" and "ICMP" (one level deeper).
$webservers = $listTr.Paragraphs(3, 1)
$webservers.InsertAfter("`rPort 80") | Out-Null
$listTr.Paragraphs(4, 1).IndentLevel = 3

$port80 = $listTr.Paragraphs(4, 1)
$port80.InsertAfter("`rICMP") | Out-Null
$listTr.Paragraphs(5, 1).IndentLevel = 3

# After that, add "Application Load Balancer" (back at the Webservers level).
$icmp1 = $listTr.Paragraphs(5, 1)
$icmp1.InsertAfter("`rApplication Load Balancer") | Out-Null
$listTr.Paragraphs(6, 1).IndentLevel = 2

# Then its two children, "HTTP" and "HTTPS".
$alb = $listTr.Paragraphs(6, 1)
$alb.InsertAfter("`rHTTP") | Out-Null
$listTr.Paragraphs(7, 1).IndentLevel = 3

$http = $listTr.Paragraphs(7, 1)
$http.InsertAfter("`rHTTPS") | Out-Null
$listTr.Paragraphs(8, 1).IndentLevel = 3

# Finally, the old "EC2" bullet becomes a third "ICMP" bullet under the ALB.
$ec2Para = $listTr.Paragraphs(9, 1)
$ec2Para.Text = "ICMP"
$ec2Para.IndentLevel = 3
